$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 37 — this shifts the existing rows 37:141 down to 38:142,
# exactly mirroring the "everything shifts down by one" pattern seen in the diff.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new weekly price record (2022-02-15).
$ws.Cells.Item(37, 1).Value = 11
$ws.Cells.Item(37, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(37, 3).Value = "Bíobío"
$ws.Cells.Item(37, 4).Value = 44607
$ws.Cells.Item(37, 5).Value = 8
$ws.Cells.Item(37, 6).Value = 100112003
$ws.Cells.Item(37, 7).Value = "Ajo"
$ws.Cells.Item(37, 8).Value = "Chino"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 400
$ws.Cells.Item(37, 11).Value = 17000
$ws.Cells.Item(37, 12).Value = 18000
$ws.Cells.Item(37, 13).Value = 17500
$ws.Cells.Item(37, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(37, 15).Value = "China"
$ws.Cells.Item(37, 16).Value = 1750
$ws.Cells.Item(37, 17).Value = 10
$ws.Cells.Item(37, 18).Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of column D.
$ws.Range("D37").NumberFormat = $ws.Range("D38").NumberFormat
